$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.224.66'
$ws.Range('E2').Value = '  -0.18%  '
$ws.Range('D3').Value = '1.866.56'
$ws.Range('E3').Value = '  -1.11%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '235.27'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.39%  '
$ws.Range('E6').Value = '  +0.04%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4664'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.47%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2831'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06543'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.01%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.29'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +6.01%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07877'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.16%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '97.55'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.82%  '
$ws.Range('D13').Value = '1.866.95'
$ws.Range('E13').Value = '  -1.15%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.107'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.53%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6731'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.81%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '280.44'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.94%  '
$ws.Range('D17').Value = '30.225.66'
$ws.Range('E17').Value = '  -0.19%  '
$ws.Range('E18').Value = '  +0.01%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.504'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.95%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.67'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.23%  '
$ws.Range('D21').Value = '2.114.73'
$ws.Range('E21').Value = '  -1.24%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.000007287'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.36%  '
$ws.Range('E23').Value = '  -0.02%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.171'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.57%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.200'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.86%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '164.86'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.67%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.15'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.90%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.924'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -4.01%  '
$ws.Range('E29').Value = '  -0.19%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.09696'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.44%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.419'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.22%  '
$ws.Range('E32').Value = '  -0.66%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.095'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.89%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.04697'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.52%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.118'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.70%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7059'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.07%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.727'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.67%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01851'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.91%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.536'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.42%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.252'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -6.20%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '73.91'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.50%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.948'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.58%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8458'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.09%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.4172'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.14%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '103.95'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.35%  '
$ws.Range('E46').Value = '  +0.06%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.184'
$ws.Range('D47').Style = 'Normal'
$ws.Range('B48').Value = 'Maker'
$ws.Range('C48').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '932.90'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -5.71%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.134'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.65%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '34.07'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.37%  '
$ws.Range('E51').Value = '  -3.40%  '
